$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:F2")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "05/08/0888"
$ws.Range("B2").Value = "5000.00"
$ws.Range("C2").Value = "5000.00"
$ws.Range("D2").Value = "5000.00"
$ws.Range("E2").Value = "5000.00"
$ws.Range("F2").Value = "20.0"
